$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = "MZyLC487"
$ws.Range("B2").Value = 231011213
$ws.Range("C2").Value = "atgqjmj10"
$ws.Range("D2").Value = "mrD2`$Y&5"
$ws.Range("F2").Value = "sgufOLrl"
$ws.Range("G2").Value = "Ssjl"

# Row 3 updates
$ws.Range("A3").Value = "WOfVH549"
$ws.Range("B3").Value = 231011198
$ws.Range("C3").Value = "wxilozm91"
$ws.Range("D3").Value = "dV`$45!Cr"
$ws.Range("F3").Value = "fyNmpqfB"
$ws.Range("G3").Value = "nRKq"
